$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new "fiyat" (price) column header in G1, matching header style of existing headers ---
$ws.Range("F1").Copy() | Out-Null
$ws.Range("G1").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false
$ws.Range("G1").Value = "fiyat"

# --- Update rows 44-62 data that got reordered/reshuffled when the price column was populated ---
$ws.Range("A44").Value = 'F By Fabrika '
$ws.Range("B44").Value = 'Termo'
$ws.Range("C44").Value = 'Deri'
$ws.Range("D44").Value = 'Türkiye'
$ws.Range("E44").Value = 'erkek'
$ws.Range("F44").Value = 'bot'
$ws.Range("A45").Value = 'F By Fabrika '
$ws.Range("B45").Value = 'Termo'
$ws.Range("D45").Value = 'Türkiye'
$ws.Range("E45").Value = 'erkek'
$ws.Range("F45").Value = 'bot'
$ws.Range("A46").Value = 'F By Fabrika '
$ws.Range("B46").Value = 'Termo'
$ws.Range("C46").Value = 'Deri'
$ws.Range("D46").Value = 'Türkiye'
$ws.Range("E46").Value = 'erkek'
$ws.Range("F46").Value = 'bot'
$ws.Range("A47").Value = 'F By Fabrika '
$ws.Range("B47").Value = 'Termo'
$ws.Range("C47").Value = 'Deri'
$ws.Range("D47").Value = 'Türkiye'
$ws.Range("E47").Value = 'erkek'
$ws.Range("A48").Value = 'F By Fabrika '
$ws.Range("B48").Value = $null
$ws.Range("D48").Value = 'Türkiye'
$ws.Range("E48").Value = 'erkek'
$ws.Range("F48").Value = 'bot'
$ws.Range("A50").Value = 'Harley Davidson '
$ws.Range("B50").Value = 'Kauçuk'
$ws.Range("D50").Value = $null
$ws.Range("E50").Value = 'kadın'
$ws.Range("F50").Value = $null
$ws.Range("A51").Value = 'Harley Davidson '
$ws.Range("B51").Value = 'Kauçuk'
$ws.Range("C51").Value = $null
$ws.Range("D51").Value = $null
$ws.Range("E51").Value = 'kadın'
$ws.Range("F51").Value = $null
$ws.Range("A52").Value = 'Harley Davidson '
$ws.Range("B52").Value = 'Kauçuk'
$ws.Range("C52").Value = $null
$ws.Range("D52").Value = $null
$ws.Range("E52").Value = 'kadın'
$ws.Range("F52").Value = $null
$ws.Range("A53").Value = 'Harley Davidson '
$ws.Range("B53").Value = 'Kauçuk'
$ws.Range("D53").Value = $null
$ws.Range("E53").Value = 'kadın'
$ws.Range("A54").Value = 'Harley Davidson '
$ws.Range("B54").Value = 'Kauçuk'
$ws.Range("C54").Value = $null
$ws.Range("D54").Value = $null
$ws.Range("E54").Value = 'kadın'
$ws.Range("F54").Value = $null
$ws.Range("A55").Value = 'Skechers '
$ws.Range("B55").Value = 'Poliüretan'
$ws.Range("D55").Value = $null
$ws.Range("E55").Value = 'erkek'
$ws.Range("F55").Value = $null
$ws.Range("A56").Value = 'Skechers '
$ws.Range("B56").Value = 'Poliüretan'
$ws.Range("D56").Value = $null
$ws.Range("E56").Value = 'erkek'
$ws.Range("F56").Value = $null
$ws.Range("A57").Value = 'Harley Davidson '
$ws.Range("B57").Value = 'Kauçuk'
$ws.Range("D57").Value = 'Türkiye'
$ws.Range("E57").Value = 'kadın'
$ws.Range("F57").Value = 'bot'
$ws.Range("A58").Value = 'Harley Davidson '
$ws.Range("B58").Value = 'Kauçuk'
$ws.Range("D58").Value = 'Türkiye'
$ws.Range("E58").Value = 'kadın'
$ws.Range("F58").Value = 'bot'
$ws.Range("A59").Value = 'New Balance '
$ws.Range("B59").Value = 'Kauçuk'
$ws.Range("C59").Value = $null
$ws.Range("D59").Value = $null
$ws.Range("F59").Value = $null
$ws.Range("A60").Value = 'New Balance '
$ws.Range("B60").Value = 'Kauçuk'
$ws.Range("C60").Value = $null
$ws.Range("D60").Value = $null
$ws.Range("F60").Value = $null
$ws.Range("A61").Value = 'Pierre Cardin '
$ws.Range("B61").Value = 'EVA'
$ws.Range("C61").Value = 'Tekstil'
$ws.Range("D61").Value = 'Türkiye'
$ws.Range("F61").Value = 'bot'
$ws.Range("A62").Value = 'Pierre Cardin '
$ws.Range("B62").Value = 'EVA'
$ws.Range("C62").Value = 'Tekstil'
$ws.Range("D62").Value = 'Türkiye'
$ws.Range("F62").Value = 'bot'

# --- Populate the new "fiyat" price values (stored as text, matching source formatting) ---
$cell = $ws.Range("G2")
$cell.NumberFormat = "@"
$cell.Value = '2799 '
$cell = $ws.Range("G3")
$cell.NumberFormat = "@"
$cell.Value = '2799 '
$cell = $ws.Range("G4")
$cell.NumberFormat = "@"
$cell.Value = '2799 '
$cell = $ws.Range("G5")
$cell.NumberFormat = "@"
$cell.Value = '2799 '
$cell = $ws.Range("G6")
$cell.NumberFormat = "@"
$cell.Value = '2799 '
$cell = $ws.Range("G7")
$cell.NumberFormat = "@"
$cell.Value = '2799 '
$cell = $ws.Range("G8")
$cell.NumberFormat = "@"
$cell.Value = '3699 '
$cell = $ws.Range("G9")
$cell.NumberFormat = "@"
$cell.Value = '2799 '
$cell = $ws.Range("G10")
$cell.NumberFormat = "@"
$cell.Value = '3699 '
$cell = $ws.Range("G11")
$cell.NumberFormat = "@"
$cell.Value = '3699 '
$cell = $ws.Range("G12")
$cell.NumberFormat = "@"
$cell.Value = '3699 '
$cell = $ws.Range("G13")
$cell.NumberFormat = "@"
$cell.Value = '3699 '
$cell = $ws.Range("G14")
$cell.NumberFormat = "@"
$cell.Value = '2799 '
$cell = $ws.Range("G15")
$cell.NumberFormat = "@"
$cell.Value = '3699 '
$cell = $ws.Range("G16")
$cell.NumberFormat = "@"
$cell.Value = '2799 '
$cell = $ws.Range("G17")
$cell.NumberFormat = "@"
$cell.Value = '2799 '
$cell = $ws.Range("G18")
$cell.NumberFormat = "@"
$cell.Value = '2799 '
$cell = $ws.Range("G19")
$cell.NumberFormat = "@"
$cell.Value = '2799 '
$cell = $ws.Range("G20")
$cell.NumberFormat = "@"
$cell.Value = '2799 '
$cell = $ws.Range("G21")
$cell.NumberFormat = "@"
$cell.Value = '2799 '
$cell = $ws.Range("G22")
$cell.NumberFormat = "@"
$cell.Value = '263912 '
$cell = $ws.Range("G23")
$cell.NumberFormat = "@"
$cell.Value = '263912 '
$cell = $ws.Range("G24")
$cell.NumberFormat = "@"
$cell.Value = '3499 '
$cell = $ws.Range("G25")
$cell.NumberFormat = "@"
$cell.Value = '3499 '
$cell = $ws.Range("G26")
$cell.NumberFormat = "@"
$cell.Value = '2499 '
$cell = $ws.Range("G27")
$cell.NumberFormat = "@"
$cell.Value = '2499 '
$cell = $ws.Range("G28")
$cell.NumberFormat = "@"
$cell.Value = '1849 '
$cell = $ws.Range("G29")
$cell.NumberFormat = "@"
$cell.Value = '1849 '
$cell = $ws.Range("G30")
$cell.NumberFormat = "@"
$cell.Value = '2449 '
$cell = $ws.Range("G31")
$cell.NumberFormat = "@"
$cell.Value = '2449 '
$cell = $ws.Range("G32")
$cell.NumberFormat = "@"
$cell.Value = '2549 '
$cell = $ws.Range("G33")
$cell.NumberFormat = "@"
$cell.Value = '2549 '
$cell = $ws.Range("G34")
$cell.NumberFormat = "@"
$cell.Value = '3749 '
$cell = $ws.Range("G35")
$cell.NumberFormat = "@"
$cell.Value = '3749 '
$cell = $ws.Range("G36")
$cell.NumberFormat = "@"
$cell.Value = '2799 '
$cell = $ws.Range("G37")
$cell.NumberFormat = "@"
$cell.Value = '2799 '
$cell = $ws.Range("G38")
$cell.NumberFormat = "@"
$cell.Value = '1699 '
$cell = $ws.Range("G39")
$cell.NumberFormat = "@"
$cell.Value = '1699 '
$cell = $ws.Range("G40")
$cell.NumberFormat = "@"
$cell.Value = '1899 '
$cell = $ws.Range("G41")
$cell.NumberFormat = "@"
$cell.Value = '1699 '
$cell = $ws.Range("G42")
$cell.NumberFormat = "@"
$cell.Value = '2999 '
$cell = $ws.Range("G43")
$cell.NumberFormat = "@"
$cell.Value = '2999 '
$cell = $ws.Range("G44")
$cell.NumberFormat = "@"
$cell.Value = '48499 '
$cell = $ws.Range("G45")
$cell.NumberFormat = "@"
$cell.Value = '60999 '
$cell = $ws.Range("G46")
$cell.NumberFormat = "@"
$cell.Value = '48499 '
$cell = $ws.Range("G47")
$cell.NumberFormat = "@"
$cell.Value = '57499 '
$cell = $ws.Range("G48")
$cell.NumberFormat = "@"
$cell.Value = '50999 '
$cell = $ws.Range("G49")
$cell.NumberFormat = "@"
$cell.Value = '48499 '
$cell = $ws.Range("G50")
$cell.NumberFormat = "@"
$cell.Value = '2999 '
$cell = $ws.Range("G51")
$cell.NumberFormat = "@"
$cell.Value = '2999 '
$cell = $ws.Range("G52")
$cell.NumberFormat = "@"
$cell.Value = '2999 '
$cell = $ws.Range("G53")
$cell.NumberFormat = "@"
$cell.Value = '3999 '
$cell = $ws.Range("G54")
$cell.NumberFormat = "@"
$cell.Value = '2999 '
$cell = $ws.Range("G55")
$cell.NumberFormat = "@"
$cell.Value = '2849 '
$cell = $ws.Range("G56")
$cell.NumberFormat = "@"
$cell.Value = '2849 '
$cell = $ws.Range("G57")
$cell.NumberFormat = "@"
$cell.Value = '3999 '
$cell = $ws.Range("G58")
$cell.NumberFormat = "@"
$cell.Value = '3999 '
$cell = $ws.Range("G59")
$cell.NumberFormat = "@"
$cell.Value = '1799 '
$cell = $ws.Range("G60")
$cell.NumberFormat = "@"
$cell.Value = '1799 '
$cell = $ws.Range("G61")
$cell.NumberFormat = "@"
$cell.Value = '89990 '
$cell = $ws.Range("G62")
$cell.NumberFormat = "@"
$cell.Value = '89990 '
